$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "grade" value for the second student row (fatima),
# matching the value already present in row 2 (Muhammad Ahmed) -> "Grade 2"
$ws.Range("I3").Value = "Grade 2"

# Move the active selection to I4 (matches where the cursor ends up after
# entering data in I3)
$ws.Range("I4").Select()
